$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: lower-case the first two headers, add four new columns ---
$ws.Range("A1").Value = "symbol"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "roll_offset"
$ws.Range("D1").Value = "roll_cycle"
$ws.Range("E1").Value = "commision"
$ws.Range("F1").Value = "multiplier"

# Carry the header style (bold font + border + centered) over to the new
# header cells so C1:F1 look like A1:B1.
$ws.Range("B1").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)

# --- Data rows: columns A/B already hold the right symbol/name text,
#     only need to populate the new C/D/E/F columns. ---

# Row 2: TA / PTA
$ws.Range("C2").Value = -220
$ws.Range("D2").Value = "1,5,9"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 5

# Row 3: M / 豆粕
$ws.Range("C3").Value = -220
$ws.Range("D3").Value = "1,5,9"
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 10

# Row 4: RB / 螺纹钢
$ws.Range("C4").Value = -220
$ws.Range("D4").Value = "1,5,10"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 10

# Row 5: FG / 玻璃
$ws.Range("C5").Value = -220
$ws.Range("D5").Value = "1,5,9"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 20

# Row 6: MA / 甲醇
$ws.Range("C6").Value = -220
$ws.Range("D6").Value = "1,5,9"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 10

# Row 7: HC / 热卷
$ws.Range("C7").Value = -220
$ws.Range("D7").Value = "1,5,10"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 10

# Row 8: V / PVC
$ws.Range("C8").Value = -220
$ws.Range("D8").Value = "1,5,9"
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 5

# Row 9: RM / 菜粕
$ws.Range("C9").Value = -220
$ws.Range("D9").Value = "1,5,9"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 10

# Row 10: AG / 沪银
$ws.Range("C10").Value = -220
$ws.Range("D10").Value = "1,6,12"
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 15

# --- Column widths for the new columns (characters) ---
$ws.Columns.Item(1).ColumnWidth = 13.36
$ws.Columns.Item(3).ColumnWidth = 19.36
$ws.Columns.Item(4).ColumnWidth = 16.08
$ws.Columns.Item(5).ColumnWidth = 21.08
$ws.Columns.Item(6).ColumnWidth = 20.08

# --- Selection matches the commit's final cursor position ---
$ws.Range("H10").Select()
